$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$samplesQuery = "MATCH (ss:study_subject)`nWITH COLLECT(ss.study_subject_id) AS all_subjects`nMATCH (samp:sample)`nMATCH (samp)-[:sample_of_study_subject]->(ss)`nMATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)`nMATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)`nMATCH (ss)<-[:diagnosis_of_study_subject]-(d)`nMATCH (d)<-[:tp_of_diagnosis]-(tp)`nWHERE s.study_acronym IN [`"A`"]  `nWITH`n    distinct lp,`n    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,`n    collect(distinct f.file_id) AS files,`n    samp, ss, s, p, all_subjects`nRETURN`n samp.sample_id AS ``Sample ID``,`n            ss.study_subject_id AS ``Case ID``,`n            p.program_acronym AS ``Program Code``,`n            s.study_acronym AS ``Arm``,`n            ss.disease_subtype AS ``Diagnosis``,`n            samp.tissue_type AS ``Tissue Type``,`n            samp.composition AS ``Tissue Composition``,`n            samp.sample_anatomic_site AS ``Sample Anatomic Site``,`n            samp.method_of_sample_procurement AS ``Sample Procurement Method``"
$filesQuery = "MATCH (f:file)-->(parent)`nMATCH (f)-[:file_of_sample]->(samp)`nMATCH (samp)-[:sample_of_study_subject]->(ss)`nMATCH (ss)-[:study_subject_of_study]->(s)`nMATCH (s)-[:study_of_program]->(p)`nMATCH (d)-[:diagnosis_of_study_subject]->(ss)`nMATCH (tp)-[:tp_of_diagnosis]->(d)`nWHERE s.study_acronym IN [`"A`"]  `nWITH`n        f, parent,p, ss, d,tp, s, samp,`n        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,`n        toInteger(floor(log(f.file_size)/log(1024))) as i,`n        2 as precision`nWITH`n        f, parent,p, ss, d,tp, s, samp,`n        f.file_size /(1024^i) AS value,`n        10^precision AS factor,`n        units[i] as unit`nWITH`n        f, parent,p, ss, d,tp, s, samp, unit,`n        round(factor * value)/factor AS size`nRETURN Distinct`n    f.file_name AS ``File Name``,`n    head(labels(samp)) AS ``Association``,`n    f.file_description AS ``Description``,`n    f.file_format AS ``File Format``,`n     CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,`n    p.program_acronym AS ``Program Code``,`n    s.study_acronym AS ``Arm``,`n    ss.study_subject_id AS ``Case ID``,`n    samp.sample_id AS ``Sample ID```n    order by f.file_name"

# Column A values first (TabName), then column B (query) -- matches shared-string
# allocation order observed in the target workbook.
$ws.Range("A3").Value2 = "SamplesTab"
$ws.Range("A4").Value2 = "FilesTab"

$ws.Range("B3").Value2 = $samplesQuery
$ws.Range("B4").Value2 = $filesQuery

$ws.Range("C3").Value2 = $ws.Range("C2").Value2
$ws.Range("C4").Value2 = $ws.Range("C2").Value2
$ws.Range("B3:C4").WrapText = $true

$ws.Range("D3").Value2 = $ws.Range("D2").Value2
$ws.Range("D4").Value2 = $ws.Range("D2").Value2

$ws.Range("E3").Value2 = $ws.Range("E2").Value2
$ws.Range("E4").Value2 = $ws.Range("E2").Value2

# Row heights (points) to accommodate the wrapped query text.
$ws.Rows.Item(2).RowHeight = 316.8
$ws.Rows.Item(3).RowHeight = 345.6
$ws.Rows.Item(4).RowHeight = 409.6

# Resize columns to fit the new content.
$ws.Columns.Item(1).ColumnWidth = 11.2504
$ws.Columns.Item(2).ColumnWidth = 73.5837
$ws.Columns.Item(3).ColumnWidth = 44.2504
$ws.Columns.Item(4).ColumnWidth = 37.417
$ws.Columns.Item(5).ColumnWidth = 36.2507

$excel.ActiveWindow.Zoom = 85
$ws.Range("C2:C4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
